$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(18).ColumnWidth = 15.4
Write-Host ($ws.Columns.Item(18).ColumnWidth)
